{"js": "// The Jinja template \"tags\" that used to live on their own separate lines\n// are combined into single paragraphs (a \"create_form\" style reflow):\n//   - \"{% set myList = [...] %}\" + \"{% for item in myList %} {% if item %}\"\n//   - \"{{ item }}\" + \"{% endif %} {% endfor %}\"\n//   - \"{% if list_1 %}\" + \"List_1 is displayed \u2013 if statement\" + \"{% endif %}\"\n//   - \"{% if list_2 %}\" + \"List_2 is displayed \u2013 if statement\" + \"{% endif %}\"\n// and one extra empty paragraph is appended at the end of the document.\n\nconst body = context.document.body;\n\n// Merge `count` consecutive paragraphs (starting with the paragraph whose\n// text matches `searchText`) into one paragraph, joining with a single\n// space and deleting the paragraphs that get folded in.\nasync function mergeParas(searchText, count) {\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  const items = paras.items;\n  let idx = -1;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === searchText) {\n      idx = i;\n      break;\n    }\n  }\n  if (idx === -1) {\n    throw new Error(\"Paragraph not found: \" + searchText);\n  }\n\n  const first = items[idx];\n  for (let k = 1; k < count; k++) {\n    const next = items[idx + k];\n    first.insertText(\" \" + next.text, \"End\");\n    next.delete();\n  }\n  await context.sync();\n}\n\n// 1) \"{% set myList = ... %}\" + \"{% for item in myList %} {% if item %}\"\nawait mergeParas(\"{% set myList = [list_1, list_2, list_3] %}\", 2);\n\n// 2) \"{{ item }}\" + \"{% endif %} {% endfor %}\"\nawait mergeParas(\"{{ item }}\", 2);\n\n// 3) \"{% if list_1 %}\" + \"List_1 is displayed \u2013 if statement\" + \"{% endif %}\"\nawait mergeParas(\"{% if list_1 %}\", 3);\n\n// 4) \"{% if list_2 %}\" + \"List_2 is displayed \u2013 if statement\" + \"{% endif %}\"\nawait mergeParas(\"{% if list_2 %}\", 3);\n\n// 5) Append one additional empty paragraph at the very end of the document.\nbody.insertParagraph(\"\", \"End\");\nawait context.sync();\n", "ps1": "# The Jinja template \"tags\" that used to live on their own separate lines\n# are combined into single paragraphs (a \"create_form\" style reflow):\n#   - \"{% set myList = [...] %}\" + \"{% for item in myList %} {% if item %}\"\n#   - \"{{ item }}\" + \"{% endif %} {% endfor %}\"\n#   - \"{% if list_1 %}\" + \"List_1 is displayed \u2013 if statement\" + \"{% endif %}\"\n#   - \"{% if list_2 %}\" + \"List_2 is displayed \u2013 if statement\" + \"{% endif %}\"\n# and one extra empty paragraph is appended at the end of the document.\n\n$d = $word.ActiveDocument\n\nfunction Merge-Paragraphs($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) \"{% set myList = ... %}\" + \"{% for item in myList %} {% if item %}\"\nMerge-Paragraphs \"{% set myList = [list_1, list_2, list_3] %}^p{% for item in myList %} {% if item %}\" \"{% set myList = [list_1, list_2, list_3] %} {% for item in myList %} {% if item %}\"\n\n# 2) \"{{ item }}\" + \"{% endif %} {% endfor %}\"\nMerge-Paragraphs \"{{ item }}^p{% endif %} {% endfor %}\" \"{{ item }} {% endif %} {% endfor %}\"\n\n# 3) \"{% if list_1 %}\" + \"List_1 is displayed \u2013 if statement\" + \"{% endif %}\"\nMerge-Paragraphs \"{% if list_1 %}^pList_1 is displayed \u2013 if statement^p{% endif %}\" \"{% if list_1 %} List_1 is displayed \u2013 if statement {% endif %}\"\n\n# 4) \"{% if list_2 %}\" + \"List_2 is displayed \u2013 if statement\" + \"{% endif %}\"\nMerge-Paragraphs \"{% if list_2 %}^pList_2 is displayed \u2013 if statement^p{% endif %}\" \"{% if list_2 %} List_2 is displayed \u2013 if statement {% endif %}\"\n\n# 5) Append one additional empty paragraph at the very end of the document.\n$d.Content.InsertParagraphAfter()\n"}
